$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 130; existing rows 130-144 shift down to 131-145.
$ws.Rows.Item(130).Insert()

# Populate the new row 130 with the weekly price-report entry (Granada, Vega Modelo de Temuco).
$ws.Range("A130").Value = 10
$ws.Range("B130").Value = "Vega Modelo de Temuco"
$ws.Range("C130").Value = "La Araucanía"
$ws.Range("D130").Value = 44769
$ws.Range("E130").Value = 9
$ws.Range("F130").Value = "Fruta"
$ws.Range("G130").Value = 100104
$ws.Range("H130").Value = "Frutos de pepita"
$ws.Range("I130").Value = 100104001
$ws.Range("J130").Value = "Granada"
$ws.Range("K130").Value = "Wonderfull"
$ws.Range("L130").Value = "Primera"
$ws.Range("M130").Value = 25
$ws.Range("N130").Value = 14000
$ws.Range("O130").Value = 14000
$ws.Range("P130").Value = 14000
$ws.Range("Q130").Value = "$/bandeja 10 kilos granel"
$ws.Range("R130").Value = "Provincia de Limarí"
$ws.Range("S130").Value = 1400
$ws.Range("T130").Value = 10
